# Risolto bug HMR e aggiunti altri log per il debug
#
# The underlying numbers for "hook" locator method (row 6, tipo "relative")
# and for "robula" locator method (row 13, tipo "relative") were wrong:
# the "Fallimenti per Fragilita" count (column E) needs correcting, which
# in turn ripples into the dependent formulas (D, G columns and the
# totals in row 20/21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the raw "Fallimenti per Fragilita" inputs.
$ws.Range("E6").Value = 1
$ws.Range("E13").Value = 0

# Force recalculation so dependent formulas (D6, G6, D13, G13, B20, D20,
# B21, D21) pick up the corrected inputs.
$excel.Calculate()

# Reflect where the user was last working in the sheet.
$ws.Range("E13").Select()

$wb.Save()
